$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new error code row data (ErrorCode / ErrorMessage for CTN0002)
$ws.Range("B12").Value = "CTN0002"
$ws.Range("C12").Value = "There is a product content already exist."

# Update the current selection / view (matches the author's saved cursor position)
$ws.Range("C13").Select() | Out-Null

# Set page setup (paper size = A4, orientation = portrait)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1 | Out-Null
